$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header values one column to the right (B1..H1 -> C1..I1)
$ws.Range("I1").Value2 = $ws.Range("H1").Value2
$ws.Range("H1").Value2 = $ws.Range("G1").Value2
$ws.Range("G1").Value2 = $ws.Range("F1").Value2
$ws.Range("F1").Value2 = $ws.Range("E1").Value2
$ws.Range("E1").Value2 = $ws.Range("D1").Value2
$ws.Range("D1").Value2 = $ws.Range("C1").Value2
$ws.Range("C1").Value2 = $ws.Range("B1").Value2

# New header for the inserted column
$ws.Range("B1").Value2 = "Name"

# Update selection to match target state
$ws.Range("B2").Select()
